$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update the "Approved/Rejected" status for row 2 from "Rejected" to "Approved"
$ws.Range("I2").Value = "Approved"

# Clear the "ReasonToReject" value for row 2 (previously "test")
$ws.Range("J2").ClearContents()

# Update the active selection to I6
$ws.Range("I6").Select()
